# Economic Dashboard V1 - weekly data refresh (2025-11-22)
# Updates the "FRED snapshot" columns (N = as-of date, Q:U = last five
# readings) for several series: the oldest reading rolls off the left,
# the remaining readings shift one column left, and a new reading is
# appended on the right (U). The as-of date in column N is bumped to the
# latest available observation date, and (row 28 only) the date cell
# picks up the "freshly updated" yellow highlight style already used by
# its neighbours - done here by copying the format from a cell that
# already carries that style, then overwriting the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28 (Mich NTM Inflation Exp / UMCSENT) ---------------------------
# N28 needs both a new date AND the yellow-highlight style (s 47 -> 48).
# Copy the format from N29 (already style 48) then set the real value.
$ws.Range("N29").Copy($ws.Range("N28"))
$ws.Range("N28").Value = 45931

$ws.Range("Q28").Value = 53.6
$ws.Range("R28").Value = 55.1
$ws.Range("S28").Value = 58.2
$ws.Range("T28").Value = 61.7
$ws.Range("U28").Value = 60.7

# --- Row 29 (5yr, 5yr Forward / T5YIFR) -----------------------------------
$ws.Range("N29").Value = 45982

$ws.Range("Q29").Value = 2.16
$ws.Range("R29").Value = 2.14
$ws.Range("S29").Value = 2.18
$ws.Range("T29").Value = 2.18
$ws.Range("U29").Value = 2.19

# --- Row 30 (10yr TIPS / T10YIE) ------------------------------------------
$ws.Range("N30").Value = 45982

$ws.Range("Q30").Value = 2.24
$ws.Range("R30").Value = 2.24
$ws.Range("S30").Value = 2.27
$ws.Range("T30").Value = 2.27
$ws.Range("U30").Value = 2.28

# --- Row 47 (FFR / DFF) ----------------------------------------------------
$ws.Range("N47").Value = 45981

# --- Row 48 (2y UST / DGS2) ------------------------------------------------
$ws.Range("N48").Value = 45981

$ws.Range("Q48").Value = 3.55
$ws.Range("R48").Value = 3.58
$ws.Range("S48").Value = 3.58
$ws.Range("T48").Value = 3.6

# --- Row 49 (5y UST / DGS5) ------------------------------------------------
$ws.Range("N49").Value = 45981

$ws.Range("Q49").Value = 3.68
$ws.Range("R49").Value = 3.71
$ws.Range("S49").Value = 3.7
$ws.Range("T49").Value = 3.72

# --- Row 50 (10y UST / DGS10) ----------------------------------------------
$ws.Range("N50").Value = 45981

$ws.Range("Q50").Value = 4.1
$ws.Range("R50").Value = 4.13
$ws.Range("S50").Value = 4.12
$ws.Range("T50").Value = 4.13

# --- Row 52 (BAA / DBAA) ----------------------------------------------------
$ws.Range("N52").Value = 45981

$ws.Range("Q52").Value = 5.9
$ws.Range("R52").Value = 5.92
$ws.Range("S52").Value = 5.91
$ws.Range("T52").Value = 5.9
